# Add "2022-Q3" data to the workbook:
#  1. Insert a new "2022-Q3" sheet (copy of "2022-Q2" format) right after "总计" / before "2022-Q2"
#     and fill it with the new quarter's fund-holding data.
#  2. Insert a new row into the "总计" (summary) sheet for "2022-Q3", shifting existing rows down
#     and renumbering the running index in column A.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. New "2022-Q3" worksheet with fund-holding detail data
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q2")
$template.Copy($template, $null)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# The template ("2022-Q2") has two data rows; the new quarter only has one.
$newSheet.Rows.Item(3).Delete()

$newSheet.Cells.Item(2, 1).Value = 0
$newSheet.Cells.Item(2, 2).NumberFormat = "@"
$newSheet.Cells.Item(2, 2).Value = "512290"
$newSheet.Cells.Item(2, 3).Value = "国泰中证生物医药ETF"
$newSheet.Cells.Item(2, 4).NumberFormat = "@"
$newSheet.Cells.Item(2, 4).Value = "40.70"
$newSheet.Cells.Item(2, 5).NumberFormat = "@"
$newSheet.Cells.Item(2, 5).Value = "99.74"
$newSheet.Cells.Item(2, 6).NumberFormat = "@"
$newSheet.Cells.Item(2, 6).Value = "2.31"
$newSheet.Cells.Item(2, 7).NumberFormat = "@"
$newSheet.Cells.Item(2, 7).Value = "0.9402"
$newSheet.Cells.Item(2, 8).Value = 7

# ---------------------------------------------------------------------------
# 2. "总计" summary sheet: insert a row for 2022-Q3 and shift the rest down
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()
$summary.Range("B2:D2").ClearFormats()

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 1
$summary.Cells.Item(2, 4).Value = 0.9399999999999999

# Match column-A's style to the rest of the index column
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

# Renumber the running index (column A) for the rows that shifted down
$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(5, 1).Value = 3
$summary.Cells.Item(6, 1).Value = 4
$summary.Cells.Item(7, 1).Value = 5
$summary.Cells.Item(8, 1).Value = 6

# Restore the originally-active sheet/selection (unrelated to the data edit itself)
$wb.Worksheets.Item("2020-Q4").Activate()
